# GUI: Stats updates JobStatus WIP
#
# The "Delete" suite (row 2) moves from JobStatus "Complete" to "Ready to
# write", so its old "Suite to manual" note in column E no longer applies
# and is cleared. The "RemovingSourceDeletesRIP" test count (C3) grows by
# one (2 -> 3). The dependent summary formulas (G1, I1, G4) recalculate
# automatically. Finally the sheet's active-cell selection moves to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Delete"): JobStatus -> "Ready to write"; clear the now-stale
# "Suite to manual" comment entirely (not just its contents).
$ws.Range("D2").Value = "Ready to write"
$ws.Range("E2").Clear()

# Row 3 ("RemovingSourceDeletesRIP"): test count 2 -> 3.
$ws.Range("C3").Value = 3

# Move the active selection to C4.
$ws.Range("C4").Select() | Out-Null
